# "Fruta / hortaliza, semanal" — weekly refresh of the Pepino ensalada sheet.
#
# The market keeps one row per day; a new day's figures (row 279) are
# inserted at the top of this variety's data block, pushing the existing
# rows 279-335 down by one (so the last existing row becomes row 336).
# No other rows on the sheet are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 279-335 down to 280-336, freeing up row 279 for the new record.
# Insert() carries the existing row formatting (e.g. the date style in
# column D) down with the shifted rows, and also the sheet's dimension
# grows to R336 automatically.
$ws.Rows.Item(279).Insert()

# New day's data for "Pepino ensalada" / Vega Modelo de Temuco.
$ws.Cells.Item(279, 1).Value = 10
$ws.Cells.Item(279, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(279, 3).Value = "La Araucanía"
$ws.Cells.Item(279, 4).Value = 44504
$ws.Cells.Item(279, 5).Value = 9
$ws.Cells.Item(279, 6).Value = 100112043
$ws.Cells.Item(279, 7).Value = "Pepino ensalada"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 235
$ws.Cells.Item(279, 11).Value = 10000
$ws.Cells.Item(279, 12).Value = 12000
$ws.Cells.Item(279, 13).Value = 10936
$ws.Cells.Item(279, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(279, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(279, 16).Value = 182
$ws.Cells.Item(279, 17).Value = 60
$ws.Cells.Item(279, 18).Value = "Hortaliza"
